# Update gh-pages to output generated at 456a3b4
# The first 3 events (南昌·SuperComic, 宜春·COMIC WORLD, 南昌·赛马娘ONLY) are no longer
# upcoming and are removed from the top of the data table; every remaining row shifts
# up by 3. A handful of "want to go" counts (column F) were refreshed to newer values,
# and the first remaining event's minimum price (G2) switched from a number to the
# text "不可售" (not for sale).

$wb = $excel.ActiveWorkbook

$sheetNames = @("展览", "全部类型")

foreach ($sheetName in $sheetNames) {
    $ws = $wb.Worksheets.Item($sheetName)

    # Remove the 3 obsolete leading data rows (rows 2, 3, 4); remaining rows shift up.
    $ws.Range("A2:I4").EntireRow.Delete() | Out-Null

    # Column A is a simple running index (row number - 1) that is independent of the
    # event content; restore it after the shift so it again reads 1..40 for rows 2..41.
    for ($r = 2; $r -le 41; $r++) {
        $ws.Cells.Item($r, 1).Value = $r - 1
    }

    # Refresh "want to go" counts (column F) that changed between scrapes.
    $ws.Range("F4").Value = 1299
    $ws.Range("F8").Value = 93
    $ws.Range("F9").Value = 9
    $ws.Range("F10").Value = 171
    $ws.Range("F11").Value = 121
    $ws.Range("F12").Value = 4378
    $ws.Range("F13").Value = 6658
    $ws.Range("F15").Value = 51
    $ws.Range("F16").Value = 92
    $ws.Range("F17").Value = 557
    $ws.Range("F19").Value = 4088
    $ws.Range("F20").Value = 445
    $ws.Range("F21").Value = 65
    $ws.Range("F23").Value = 2663
    $ws.Range("F27").Value = 334
    $ws.Range("F28").Value = 343
    $ws.Range("F30").Value = 211
    $ws.Range("F32").Value = 1605
    $ws.Range("F33").Value = 1008
    $ws.Range("F35").Value = 116
    $ws.Range("F36").Value = 74
    $ws.Range("F37").Value = 523
    $ws.Range("F41").Value = 616

    # Minimum price for the now-first row became unavailable for sale.
    $ws.Range("G2").Value = "不可售"
}
